$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-29 01:46:29"
}
